# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the f2fa4ce0-... entry (row 4) on both the zh-cn and de-de report
# sheets, matching a later handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-22 09:12:54"
$wsZhCn.Range("G4").Value = "2016-02-22 09:13:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-22 09:13:09"
$wsDeDe.Range("G4").Value = "2016-02-22 09:14:15"
